$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the footer block (underline + signature labels) down one row by inserting
# a blank row in an area with no existing cell data, so no new styles are created.
$ws.Rows("20").Insert()

# Push the current last data row (period 2506, bottom-border style) down into row 19.
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))
$excel.CutCopyMode = 0

# Re-create row 18 using the same look as the row above it (period 2505, middle style),
# then overwrite row 18 with the worker's data for the new period.
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))
$excel.CutCopyMode = 0

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "13881811"
$ws.Range("D18").Value = "JESUS ALBERTO CEBALLOS ALVAREZ"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Update the period labels for the whole (now 4-row) table, newest period first
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"
$ws.Range("E19").Value = "2504"

# Update totals to reflect the newly added period
$ws.Range("E11").Value = 227760
$ws.Range("F13").Value = 4
